$wb = $excel.ActiveWorkbook

# "OFF" sheet - Week 13 update for row 2 ("H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 453
$wsOff.Range("C2").Value = 317
$wsOff.Range("D2").Value = 95
$wsOff.Range("E2").Value = 37

# "DEF" sheet - Week 13 update for row 2 ("H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 385
$wsDef.Range("C2").Value = 284
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 49
$wsDef.Range("F2").Value = 5
